# Fruta / hortaliza, semanal
# Insert 6 new weekly rows of apple (Manzana) price data for the
# Vega Central Mapocho de Santiago market, pushing the existing rows
# 1137-1189 down to 1143-1195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert six blank rows right before the current row 1137.
$ws.Rows("1137:1142").Insert()

# Shared (unchanged) column values for every row in this block.
$mercadoId   = 9
$mercado     = "Vega Central Mapocho de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100104
$producto    = "Frutos de pepita"
$categoriaId = 100104002
$categoria   = "Manzana"
$fecha       = 44509

function Set-Row {
    param(
        [int]$RowNum,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($RowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($RowNum, 2).Value  = $mercado
    $ws.Cells.Item($RowNum, 3).Value  = $region
    $ws.Cells.Item($RowNum, 4).Value  = $fecha
    $ws.Cells.Item($RowNum, 5).Value  = $codreg
    $ws.Cells.Item($RowNum, 6).Value  = $tipo
    $ws.Cells.Item($RowNum, 7).Value  = $productoId
    $ws.Cells.Item($RowNum, 8).Value  = $producto
    $ws.Cells.Item($RowNum, 9).Value  = $categoriaId
    $ws.Cells.Item($RowNum, 10).Value = $categoria
    $ws.Cells.Item($RowNum, 11).Value = $Variedad
    $ws.Cells.Item($RowNum, 12).Value = $Calidad
    $ws.Cells.Item($RowNum, 13).Value = $Volumen
    $ws.Cells.Item($RowNum, 14).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 15).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 16).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value = $Unidad
    $ws.Cells.Item($RowNum, 18).Value = $Origen
    $ws.Cells.Item($RowNum, 19).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value = $KgUnidad
}

Set-Row 1137 "Fuji"         "Calibre 90" 520 15000 16000 15615 "`$/caja 18 kilos embalada" "Región de O'Higgins" 868 18
Set-Row 1138 "Fuji"         "Primera"    350 12000 12000 12000 "`$/caja 18 kilos granel"   "Provincia de Curicó"   667 18
Set-Row 1139 "Granny Smith" "Calibre 90" 470 15000 16000 15532 "`$/caja 18 kilos embalada" "Provincia de Curicó"   863 18
Set-Row 1140 "Pink Lady"    "Primera"    350 11000 11000 11000 "`$/caja 18 kilos granel"   "Región de O'Higgins" 611 18
Set-Row 1141 "Scarlett"     "Calibre 90" 300 15000 15000 15000 "`$/caja 18 kilos embalada" "Provincia de Linares"  833 18
Set-Row 1142 "Scarlett"     "Primera"    300 12000 12000 12000 "`$/caja 18 kilos granel"   "Región de O'Higgins" 667 18
